$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRow = 65

# Text columns (A-C): company, location, type
$ws.Cells.Item($newRow, 1).Value = "CompaNanny"
$ws.Cells.Item($newRow, 2).Value = "CompaNanny Rembrandtlaan"
$ws.Cells.Item($newRow, 3).Value = "BSO"

# Date column (D) is stored as a plain text string in this workbook
# (e.g. "2023-12-18"), not an Excel date serial. Force text formatting
# before assigning, then clear the formatting again so the new row's
# cell carries no leftover style index, matching the rest of the sheet.
$ws.Cells.Item($newRow, 4).NumberFormat = "@"
$ws.Cells.Item($newRow, 4).Value = "2023-12-18"
$ws.Cells.Item($newRow, 4).ClearFormats()

# Numeric score columns (E-J) are all 0 for this new entry
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 0
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 0
